# Update "want to go" counts (column F) across sheets, reflecting a
# fresh scrape snapshot ("Update gh-pages to output generated at 456a3b4").
# Only the F-column numeric values change; everything else is untouched.
#
# NOTE: named parameters (-Foo bar) on custom functions are not reliably
# bound by this runtime, so the updates are written as plain, direct
# statements instead of going through a helper function.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 848
$ws.Range("F3").Value = 574
$ws.Range("F8").Value = 64
$ws.Range("F11").Value = 1239
$ws.Range("F13").Value = 93
$ws.Range("F14").Value = 926
$ws.Range("F15").Value = 901
$ws.Range("F19").Value = 660
$ws.Range("F20").Value = 850
$ws.Range("F21").Value = 1769
$ws.Range("F22").Value = 3350
$ws.Range("F23").Value = 985
$ws.Range("F25").Value = 2378
$ws.Range("F27").Value = 27
$ws.Range("F28").Value = 3274
$ws.Range("F30").Value = 814
$ws.Range("F31").Value = 23
$ws.Range("F32").Value = 2020
$ws.Range("F33").Value = 98
$ws.Range("F34").Value = 768
$ws.Range("F36").Value = 150
$ws.Range("F37").Value = 100
$ws.Range("F38").Value = 106
$ws.Range("F39").Value = 1174
$ws.Range("F40").Value = 1846
$ws.Range("F41").Value = 442
$ws.Range("F44").Value = 221
$ws.Range("F46").Value = 204
$ws.Range("F47").Value = 62

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 103

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 175

# 全部类型 (All types - combined listing)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 848
$ws.Range("F3").Value = 574
$ws.Range("F7").Value = 1239
$ws.Range("F8").Value = 926
$ws.Range("F9").Value = 901
$ws.Range("F16").Value = 850
$ws.Range("F17").Value = 1769
$ws.Range("F18").Value = 3350
$ws.Range("F19").Value = 985
$ws.Range("F21").Value = 2378
$ws.Range("F22").Value = 27
$ws.Range("F23").Value = 3274
$ws.Range("F25").Value = 814
$ws.Range("F27").Value = 23
$ws.Range("F28").Value = 2020
$ws.Range("F32").Value = 98
$ws.Range("F33").Value = 103
$ws.Range("F34").Value = 768
$ws.Range("F36").Value = 150
$ws.Range("F37").Value = 100
$ws.Range("F38").Value = 106
$ws.Range("F41").Value = 1174
$ws.Range("F42").Value = 1846
$ws.Range("F45").Value = 442
$ws.Range("F47").Value = 221
$ws.Range("F49").Value = 204
$ws.Range("F50").Value = 62

Write-Host "Done updating F-column values."
